$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Unmerge the C24:D24 cell (it will be split into two separate param cells: p1 / p2)
$ws.Range("C24:D24").UnMerge()

# Update the description cells for the test sub-table: param1/param2 -> p1/p2
$ws.Range("C24").Value = "p1"
$ws.Range("D24").Value = "p2"

# Remove horizontal centering on these two cells (keep vertical center + wrap text),
# use General alignment (1) so the horizontal attribute is simply dropped.
$ws.Range("C24").HorizontalAlignment = 1
$ws.Range("D24").HorizontalAlignment = 1

# Update the test result value in C26: 1 -> 4
$ws.Range("C26").Value = 4

# Update the active selection to match target state
$ws.Range("X18").Select()
